# "Added sound to Assest list.xlsx"
# Insert a new "Fly passive buzzing" sound entry (with a couple of blank
# spacer rows) above the "Prefabs" section on the Knight+Player sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the "Prefabs" block (old rows 35-40) down by 3 rows, inheriting the
# formatting of the row above (matches the blank B33/B34-style rows).
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

# New sound line item; rows 36-37 are left blank as spacing before "Prefabs".
$ws.Range("B35").Value = "Fly passive buzzing"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("B36").Select()
